$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New standalone data row used by the Manage Schedule test case
$ws.Range("B19").Value = "SAMIKHAN27980"

# Update login id value for the valid-login test case, and the new
# Limit Management test case row, to the same user id
$ws.Range("B2").Value = "MOHSIN7812"
$ws.Range("B14").Value = "MOHSIN7812"

# Move the active selection like the authored workbook
$ws.Range("G24").Select()
